# Update market price / profit figures in the Hyperion_Profits sheets
# (scheduled runner refresh of currentAveragePrice / LeveProfit columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1810.375
$ws.Range("J43").Value = 1688.5
$ws.Range("L43").Value = 1688.5
$ws.Range("N43").Value = -1826.5

$ws.Range("H52").Value = 1728.3334
$ws.Range("I52").Value = 592.5
$ws.Range("K52").Value = 1777.5
$ws.Range("M52").Value = -1617.5

$ws.Range("H86").Value = 2798
$ws.Range("I86").Value = 1990
$ws.Range("K86").Value = 1990
$ws.Range("M86").Value = -867

$ws.Range("H89").Value = 2798
$ws.Range("I89").Value = 1990
$ws.Range("K89").Value = 9950
$ws.Range("M89").Value = -4334

$ws.Range("H125").Value = 3125.2075
$ws.Range("J125").Value = 3534.2195
$ws.Range("L125").Value = 31807.9755
$ws.Range("N125").Value = -36727.9755

$ws.Range("H132").Value = 2079.4
$ws.Range("I132").Value = 2082.138
$ws.Range("K132").Value = 6246.414
$ws.Range("M132").Value = -3716.414

$ws.Range("H137").Value = 2541.3914
$ws.Range("I137").Value = 2300.5652
$ws.Range("J137").Value = 2782.2173
$ws.Range("K137").Value = 6901.6956
$ws.Range("L137").Value = 8346.651899999999
$ws.Range("M137").Value = -4351.6956
$ws.Range("N137").Value = -13446.6519


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4031.137
$ws.Range("I32").Value = 2631.5132
$ws.Range("K32").Value = 2631.5132
$ws.Range("M32").Value = -2344.5132

$ws.Range("H61").Value = 1709.8077
$ws.Range("I61").Value = 1623.1666
$ws.Range("K61").Value = 1623.1666
$ws.Range("M61").Value = -1411.1666

$ws.Range("H74").Value = 68289.64999999999
$ws.Range("I74").Value = 57693.76
$ws.Range("J74").Value = 90541
$ws.Range("K74").Value = 57693.76
$ws.Range("L74").Value = 90541
$ws.Range("M74").Value = -56819.76
$ws.Range("N74").Value = -92289

$ws.Range("H77").Value = 68289.64999999999
$ws.Range("I77").Value = 57693.76
$ws.Range("J77").Value = 90541
$ws.Range("K77").Value = 288468.8
$ws.Range("L77").Value = 452705
$ws.Range("M77").Value = -284100.8
$ws.Range("N77").Value = -461441

$ws.Range("H136").Value = 1709.8077
$ws.Range("I136").Value = 1623.1666
$ws.Range("K136").Value = 4869.4998
$ws.Range("M136").Value = -2319.4998


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2666.6667
$ws.Range("I5").Value = 2000
$ws.Range("K5").Value = 2000
$ws.Range("M5").Value = -1887

$ws.Range("H94").Value = 4764.6743
$ws.Range("I94").Value = 683.1875
$ws.Range("J94").Value = 16638.092
$ws.Range("K94").Value = 683.1875
$ws.Range("L94").Value = 16638.092
$ws.Range("M94").Value = -232.1875
$ws.Range("N94").Value = -17540.092

$ws.Range("H96").Value = 5646.5
$ws.Range("I96").Value = 5646.5
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 5646.5
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -2900.5
$ws.Range("N96").ClearContents()


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2980.9512
$ws.Range("I31").Value = 1967.0667
$ws.Range("J31").Value = 5746.091
$ws.Range("K31").Value = 1967.0667
$ws.Range("L31").Value = 5746.091
$ws.Range("M31").Value = -1672.0667
$ws.Range("N31").Value = -6336.091

$ws.Range("H33").Value = 3115.4285
$ws.Range("I33").Value = 1143.25
$ws.Range("J33").Value = 5745
$ws.Range("K33").Value = 1143.25
$ws.Range("L33").Value = 5745
$ws.Range("M33").Value = -764.25
$ws.Range("N33").Value = -6503

$ws.Range("H34").Value = 2980.9512
$ws.Range("I34").Value = 1967.0667
$ws.Range("J34").Value = 5746.091
$ws.Range("K34").Value = 1967.0667
$ws.Range("L34").Value = 5746.091
$ws.Range("M34").Value = -1765.0667
$ws.Range("N34").Value = -6150.091


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 167518.33
$ws.Range("I9").Value = 500005
$ws.Range("J9").Value = 1275
$ws.Range("K9").Value = 1500015
$ws.Range("L9").Value = 3825
$ws.Range("M9").Value = -1499791
$ws.Range("N9").Value = -4273

$ws.Range("H10").Value = 13.833333
$ws.Range("J10").Value = 12.5
$ws.Range("L10").Value = 37.5
$ws.Range("N10").Value = -315.5

$ws.Range("H13").Value = 150.16667
$ws.Range("I13").Value = 7
$ws.Range("K13").Value = 21
$ws.Range("M13").Value = 147

$ws.Range("H17").Value = 1737.75
$ws.Range("I17").Value = 1169.375
$ws.Range("J17").Value = 2116.6667
$ws.Range("K17").Value = 3508.125
$ws.Range("L17").Value = 6350.000100000001
$ws.Range("M17").Value = -3339.125
$ws.Range("N17").Value = -6688.000100000001

$ws.Range("H37").Value = 42039.4
$ws.Range("J37").Value = 42039.4
$ws.Range("L37").Value = 126118.2
$ws.Range("N37").Value = -126342.2


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 15499.5
$ws.Range("I21").Value = 10999
$ws.Range("K21").Value = 10999
$ws.Range("M21").Value = -10826

$ws.Range("H30").Value = 15499.5
$ws.Range("I30").Value = 10999
$ws.Range("K30").Value = 10999
$ws.Range("M30").Value = -10894

$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -21064

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1451.2
$ws.Range("I9").Value = 500
$ws.Range("J9").Value = 1689
$ws.Range("K9").Value = 500
$ws.Range("L9").Value = 1689
$ws.Range("M9").Value = -276
$ws.Range("N9").Value = -2137

$ws.Range("H132").Value = 6643.9
$ws.Range("I132").Value = 6773.5
$ws.Range("K132").Value = 20320.5
$ws.Range("M132").Value = -17790.5

$ws.Range("H136").Value = 46202.74
$ws.Range("I136").Value = 53059.64
$ws.Range("K136").Value = 159178.92
$ws.Range("M136").Value = -156628.92


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 19500
$ws.Range("I40").Value = 21000
$ws.Range("J40").Value = 18000
$ws.Range("K40").Value = 21000
$ws.Range("L40").Value = 18000
$ws.Range("M40").Value = -20851
$ws.Range("N40").Value = -18298

$ws.Range("H80").Value = 47466.668
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 47466.668
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 47466.668
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -49462.668

$ws.Range("H83").Value = 47466.668
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 47466.668
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 142400.004
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -152384.004

$ws.Range("H104").Value = 45000
$ws.Range("J104").Value = 45000
$ws.Range("L104").Value = 45000
$ws.Range("N104").Value = -51988

$ws.Range("H141").Value = 125000
$ws.Range("J141").Value = 125000
$ws.Range("L141").Value = 125000
$ws.Range("N141").Value = -135360

